$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4119732975959778
$ws.Range("B1").Value = 2.755685806274414
$ws.Range("C1").Value = 6.085403442382812
$ws.Range("D1").Value = 1.714774012565613
$ws.Range("E1").Value = 1.018853425979614
